# daily auto push: 2026-02-09 22:56 UTC
# Insert the latest reading as a new row right before the existing
# 2026/12/29 block (row 793), pushing the rest of the log down by one
# row and growing the used range from D834 to D835.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 793

# Push rows 793:834 down to 794:835 to make room for the new entry.
$ws.Rows.Item($newRow).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a
# real date value, so force text formatting before writing it -
# otherwise a date-looking string gets auto-converted to a serial date.
# Reset the style back to Normal afterwards so the new cell ends up
# with the same (default) formatting as the rest of the column.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/02/10"
$ws.Cells.Item($newRow, 1).Style = "Normal"
$ws.Cells.Item($newRow, 2).Value = "火"
$ws.Cells.Item($newRow, 3).Value = 6
$ws.Cells.Item($newRow, 4).Value = 201
